# Align text of birthdays
# - Rename Planilha1 -> DataBase
# - Add a new "Assinatura" sheet (signature block) after DataBase
# - Fix a birthdate value on DataBase!E17
# - Update the DataBase sheet's view/selection

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Rename the original sheet and insert the new "Assinatura" sheet
# ---------------------------------------------------------------
$dataBase = $wb.Worksheets.Item(1)
$dataBase.Name = "DataBase"

$assinatura = $wb.Worksheets.Add($null, $dataBase)
$assinatura.Name = "Assinatura"

# ---------------------------------------------------------------
# 2) Fix up a birthdate value on the DataBase sheet
# ---------------------------------------------------------------
$dataBase.Range("E17").Value = 21212

# ---------------------------------------------------------------
# 3) Populate the Assinatura (signature) sheet
# ---------------------------------------------------------------
$assinatura.Range("A1").Value = "Assinatura"
$assinatura.Range("A2").Value = "Marciano de Freitas Matos"
$assinatura.Range("A3").Value = "SR Brasília sul"
$assinatura.Range("A4").Value = "marciano.matos@caixa.gov.br"
$assinatura.Range("A5").Value = "Superintendente de Rede"

# Header cell (A1): bold white text on a solid blue (Accent 1) fill,
# thin accent-colored border on right/top/bottom.
$assinatura.Range("A1").Font.Bold = $true
$assinatura.Range("A1").Font.ThemeColor = 2
$assinatura.Range("A1").Interior.ThemeColor = 5
$assinatura.Range("A1").Borders.Item(10).LineStyle = 1
$assinatura.Range("A1").Borders.Item(8).LineStyle = 1
$assinatura.Range("A1").Borders.Item(9).LineStyle = 1
$assinatura.Range("A1").Borders.Item(10).Color = 14410198
$assinatura.Range("A1").Borders.Item(8).Color = 14410198
$assinatura.Range("A1").Borders.Item(9).Color = 14410198

# Shaded rows (A2, A4): light-blue (Accent 1, lighter 80%) fill with the
# same thin border, no special font.
foreach ($addr in @("A2", "A4")) {
    $rng = $assinatura.Range($addr)
    $rng.Interior.Color = 14410198
    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(10).Color = 14410198
    $rng.Borders.Item(8).Color = 14410198
    $rng.Borders.Item(9).Color = 14410198
}

# Plain bordered row (A3): no fill, same thin border.
$assinatura.Range("A3").Borders.Item(10).LineStyle = 1
$assinatura.Range("A3").Borders.Item(8).LineStyle = 1
$assinatura.Range("A3").Borders.Item(9).LineStyle = 1
$assinatura.Range("A3").Borders.Item(10).Color = 14410198
$assinatura.Range("A3").Borders.Item(8).Color = 14410198
$assinatura.Range("A3").Borders.Item(9).Color = 14410198

$assinatura.Columns.Item(1).ColumnWidth = 28.14

$assinatura.Range("A6").Select()

# ---------------------------------------------------------------
# 4) Restore DataBase as the active sheet / selection
# ---------------------------------------------------------------
$dataBase.Select()
$dataBase.Range("E13").Select()
$excel.ActiveWindow.ScrollColumn = 3
